# Rename translation-column headers:
#   "display.text"  -> "display.title.text"  (on the "initial" and "survey" sheets)
#   "display.title" -> "display.title.text"  (on the "settings" sheet)
# The "choices" sheet keeps its "display.text" header unchanged.

$wb = $excel.ActiveWorkbook

$wsInitial  = $wb.Worksheets.Item("initial")
$wsSurvey   = $wb.Worksheets.Item("survey")
$wsSettings = $wb.Worksheets.Item("settings")

$wsInitial.Range("C1").Value  = "display.title.text"
$wsSurvey.Range("D1").Value   = "display.title.text"
$wsSettings.Range("C1").Value = "display.title.text"

# Update the view state: selections move to the cells the author was
# last looking at, and the "settings" sheet becomes the active tab.
$wsInitial.Range("C2").Select()
$wsSurvey.Range("D2").Select()

$wsSettings.Activate()
$wsSettings.Range("C2").Select()
